$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 1.0.0 -> 1.0.1
$ws.Range("B3").Value = "1.0.1"

# Contact value changed
$ws.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Insert a new "Jurisdiction" row right after "Contact" (row 11), pushing
# Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting used by the other data rows (border/alignment style).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# --- Rename the second sheet ---
$ws2 = $wb.Worksheets.Item("Include from unknown")
$ws2.Name = "Include #0"
